$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C slightly (bestFit grew with new data)
$ws.Columns.Item(3).ColumnWidth = 8.875

# Add the new trade row (row 8) - copy formats from row 7 first, to reuse existing styles
$ws.Range("A7:I7").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)

$ws.Range("A8").Value = 42649.644814814812
$ws.Range("B8").Value = $true
$ws.Range("C8").Value = 10068.83
$ws.Range("D8").Value = 9957.7999999999993
$ws.Range("E8").Value = 18.829999999999998
$ws.Range("F8").Value = 19.25
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = 2.23
$ws.Range("I8").Value = $false
